$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'71.117.65"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "`'3.547.49"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "`'0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "`'618.87"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "`'174.28"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").Value = "`'3.540.94"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "`'0.613"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").Value = "`'0.999"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "`'0.199"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").Value = "`'7.23"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "`'0.591"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "`'46.80"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "`'0.0000278"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "`'4.112.37"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "`'8.47"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "`'614.77"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "`'3.544.69"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "`'71.069.95"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "`'0.122"
$ws.Range("D21").Value = "`'17.85"
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("D22").Value = "`'0.891"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "`'9.10"
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("D24").Value = "`'15.79"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "`'98.25"
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "`'2.61"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "`'33.96"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("D30").Value = "`'9.20"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("D31").Value = "`'3.05"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "`'8.21"
$ws.Range("E32").Value = "  -3.27%  "
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").Value = "`'6.89"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "`'615.73"
$ws.Range("E35").Value = "  +7.42%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").Value = "`'0.0479"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").Value = "`'57.03"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "`'0.0₃0744"
$ws.Range("E43").Value = "  +6.15%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "`'3.382.18"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "`'2.99"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "`'0.316"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "`'32.36"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").Value = "`'2.58"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "`'133.60"
$ws.Range("E50").Value = "  -0.13%  "
